# Append 19 new "tensorflow/ranking" dataset rows (ids 139-157) to the
# "dataset" sheet, rows 140-158, mirroring the existing data shape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$repoUrl  = "https://github.com/tensorflow/ranking"
$repoName = "ranking"
$repoAuth = "tensorflow"
$startDate = "12/03/2018"

$firstNewRow = 140
$lastNewRow  = 158

for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {

    $id = $r - 1

    # id column is a real number, formatted/bordered like the rows above it
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 1).VerticalAlignment = -4160
    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 1).Borders.LineStyle = 1

    $ws.Cells.Item($r, 2).Value = $repoUrl
    $ws.Cells.Item($r, 3).Value = $repoName
    $ws.Cells.Item($r, 4).Value = $repoAuth
    $ws.Cells.Item($r, 5).Value = $startDate

    # smell flag columns (OSE BCE PDE SV OS SD RS TFS UI TC) stored as text "0"/"1"
    $ws.Cells.Item($r, 6).Value  = "'0"   # OSE
    $ws.Cells.Item($r, 7).Value  = "'0"   # BCE
    $ws.Cells.Item($r, 8).Value  = "'0"   # PDE
    $ws.Cells.Item($r, 9).Value  = "'1"   # SV

    if ($r -le 141) {
        $ws.Cells.Item($r, 10).Value = "'0"   # OS
    } else {
        $ws.Cells.Item($r, 10).Value = "'1"   # OS
    }

    $ws.Cells.Item($r, 11).Value = "'1"   # SD
    $ws.Cells.Item($r, 12).Value = "'0"   # RS
    $ws.Cells.Item($r, 13).Value = "'0"   # TFS
    $ws.Cells.Item($r, 14).Value = "'1"   # UI
    $ws.Cells.Item($r, 15).Value = "'0"   # TC
}

Write-Host "Appended rows $firstNewRow..$lastNewRow"
